$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was reported for "Poroto granado" (Vega Modelo de
# Temuco). It belongs right after the current row 32, so insert a fresh row
# at position 33 — this pushes the existing rows 33..101 down to 34..102
# (the former last data row, row 101, becomes row 102) and grows the used
# range to A1:R102, exactly like Excel's own "Insert Row" command would.
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new record's data.
$ws.Range('A33').Value = 10
$ws.Range('B33').Value = 'Vega Modelo de Temuco'
$ws.Range('C33').Value = 'La Araucanía'
$ws.Range('D33').Value = 44952
$ws.Range('E33').Value = 9
$ws.Range('F33').Value = 100112030
$ws.Range('G33').Value = 'Poroto granado'
$ws.Range('H33').Value = 'Sin especificar'
$ws.Range('I33').Value = 'Primera'
$ws.Range('J33').Value = 135
$ws.Range('K33').Value = 45000
$ws.Range('L33').Value = 45000
$ws.Range('M33').Value = 45000
$ws.Range('N33').Value = '$/saco 25 kilos'
$ws.Range('O33').Value = 'Región del Maule'
$ws.Range('P33').Value = 1800
$ws.Range('Q33').Value = 25
$ws.Range('R33').Value = 'Hortaliza'
